# Updates the cryptos sheet Price (D) / Volume(1h) (E) columns for rows 2-51
# to the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price/Volume cells are stored as plain text (e.g. "0.420", "7.20",
# "  -3.51%  "). Assigning a string that LOOKS like a plain number through
# Range.Value makes Excel silently reinterpret it as a Number (dropping
# trailing zeros, switching to scientific notation, etc). Force those cells
# to Text format before the write, then clear the format again afterwards so
# the cell style is left exactly as it was (only the content changes).
$textForced = @(
    "D5", "D6", "D8", "D9", "D14", "D15", "D17", "D19",
    "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D30",
    "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39",
    "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49",
    "D51"
)
foreach ($addr in $textForced) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "55.639.15"
$ws.Range("E2").Value = "  -3.51%  "
$ws.Range("D3").Value = "2.915.48"
$ws.Range("E3").Value = "  -3.74%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "506.44"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").Value = "133.19"
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.420"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").Value = "7.20"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("D12").Value = "3.412.17"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "25.81"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").Value = "0.0000159"
$ws.Range("E15").Value = "  -4.14%  "
$ws.Range("D16").Value = "55.609.07"
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("D17").Value = "6.01"
$ws.Range("E17").Value = "  -3.49%  "
$ws.Range("D18").Value = "2.913.51"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  -4.30%  "
$ws.Range("D21").Value = "314.84"
$ws.Range("E21").Value = "  -5.72%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "0.484"
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("D24").Value = "62.55"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("D25").Value = "3.031.14"
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  -4.34%  "
$ws.Range("D28").Value = "0.0₃0843"
$ws.Range("E28").Value = "  -9.14%  "
$ws.Range("E29").Value = "  -6.80%  "
$ws.Range("D30").Value = "6.92"
$ws.Range("E30").Value = "  -7.29%  "
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "19.73"
$ws.Range("E32").Value = "  -5.34%  "
$ws.Range("D33").Value = "1.13"
$ws.Range("E33").Value = "  -7.42%  "
$ws.Range("D34").Value = "148.79"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").Value = "4.38"
$ws.Range("E35").Value = "  -7.06%  "
$ws.Range("D36").Value = "5.60"
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("D37").Value = "24.56"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("E38").Value = "  -7.63%  "
$ws.Range("D39").Value = "0.0646"
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "36.31"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").Value = "3.70"
$ws.Range("E42").Value = "  -4.59%  "
$ws.Range("D43").Value = "0.636"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "2.106.25"
$ws.Range("E44").Value = "  -8.71%  "
$ws.Range("D45").Value = "1.33"
$ws.Range("E45").Value = "  -6.87%  "
$ws.Range("D46").Value = "5.90"
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("D47").Value = "0.915"
$ws.Range("E47").Value = "  -7.62%  "
$ws.Range("D48").Value = "0.0231"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").Value = "18.63"
$ws.Range("E49").Value = "  -5.26%  "
$ws.Range("E50").Value = "  -6.80%  "
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  -8.79%  "

foreach ($addr in $textForced) { $ws.Range($addr).ClearFormats() }
